$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.351959228515625
$ws.Range("B1").Value = 4.055415630340576
$ws.Range("C1").Value = 1.941134572029114
$ws.Range("D1").Value = 1.470391273498535
$ws.Range("E1").Value = 1.308032989501953
